$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '97.145.07'
$ws.Range("E2").Value = '  +4.80%  '
# Row 3
$ws.Range("D3").Value = '3.130.10'
$ws.Range("E3").Value = '  +0.63%  '
# Row 4
$ws.Range("E4").Value = '  +0.08%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.59'
$ws.Range("E5").Value = '  +1.68%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '612.62'
$ws.Range("E6").Value = '  -0.27%  '
# Row 7
$ws.Range("E7").Value = '  +1.74%  '
# Row 8
$ws.Range("E8").Value = '  -1.14%  '
# Row 9
$ws.Range("E9").Value = '  +0.09%  '
# Row 10
$ws.Range("D10").Value = '3.123.38'
$ws.Range("E10").Value = '  +0.54%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.782'
$ws.Range("E11").Value = '  -1.16%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.198'
$ws.Range("E12").Value = '  -0.11%  '
# Row 13
$ws.Range("D13").Value = '96.848.46'
$ws.Range("E13").Value = '  +4.83%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000241'
$ws.Range("E14").Value = '  -1.40%  '
# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.99'
$ws.Range("E15").Value = '  +0.20%  '
# Row 16
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.43'
$ws.Range("E16").Value = '  +0.11%  '
# Row 17
$ws.Range("D17").Value = '3.709.74'
$ws.Range("E17").Value = '  +0.72%  '
# Row 18
$ws.Range("D18").Value = '3.130.10'
$ws.Range("E18").Value = '  +0.67%  '
# Row 19
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.52'
$ws.Range("E19").Value = '  -7.46%  '
# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '524.62'
$ws.Range("E20").Value = '  +19.44%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.55'
$ws.Range("E21").Value = '  +0.33%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.70'
$ws.Range("E22").Value = '  -2.28%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000194'
$ws.Range("E23").Value = '  -4.96%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.86'
$ws.Range("E24").Value = '  -2.47%  '
# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.94'
$ws.Range("E25").Value = '  +4.18%  '
# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.48'
$ws.Range("E26").Value = '  -1.81%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.60'
$ws.Range("E27").Value = '  +0.18%  '
# Row 28
$ws.Range("D28").Value = '3.285.03'
$ws.Range("E28").Value = '  +0.35%  '
# Row 29
$ws.Range("E29").Value = '  +0.17%  '
# Row 30
$ws.Range("E30").Value = '  +0.29%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.175'
$ws.Range("E31").Value = '  -4.27%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.125'
$ws.Range("E32").Value = '  +0.44%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.02'
$ws.Range("E33").Value = '  -1.49%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.71'
$ws.Range("E34").Value = '  +3.99%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.865'
$ws.Range("E35").Value = '  -13.43%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.153'
$ws.Range("E36").Value = '  -8.04%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.34'
$ws.Range("E37").Value = '  -9.44%  '
# Row 38
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '483.46'
$ws.Range("E38").Value = '  +4.18%  '
# Row 39
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.87'
$ws.Range("E39").Value = '  -0.98%  '
# Row 40
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.24'
$ws.Range("E40").Value = '  +1.37%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.440'
$ws.Range("E41").Value = '  +2.67%  '
# Row 42
$ws.Range("E42").Value = '  -3.89%  '
# Row 43
$ws.Range("E43").Value = '  -10.98%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.19'
$ws.Range("E45").Value = '  -4.57%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.18'
$ws.Range("E46").Value = '  +1.50%  '
# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.700'
$ws.Range("E47").Value = '  +2.70%  '
# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.91'
$ws.Range("E48").Value = '  +4.08%  '
# Row 49
$ws.Range("E49").Value = '  +3.66%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.40'
$ws.Range("E50").Value = '  +1.40%  '
# Row 51
$ws.Range("E51").Value = '  +0.12%  '
